$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("O1").Value = "F1 train"

# Column O (Validation / F1 train) value updates for rows that only change O
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 0.9714285714285714
$ws.Range("O4").Value = 0.9859154929577465
$ws.Range("O5").Value = 1
$ws.Range("O7").Value = 0.927536231884058
$ws.Range("O8").Value = 0.9166666666666666
$ws.Range("O9").Value = 0.9859154929577465
$ws.Range("O10").Value = 1
$ws.Range("O12").Value = 0.8571428571428571
$ws.Range("O13").Value = 0.9855072463768116
$ws.Range("O14").Value = 0.8571428571428571
$ws.Range("O15").Value = 0.78125

# Row 6 (MLP, technique 5) - full metrics refresh + parameter change
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 0.65
$ws.Range("J6").Value = 0.6666666666666666
$ws.Range("K6").Value = 0.7777777777777778
$ws.Range("L6").Value = 0.5833333333333334
$ws.Range("M6").Value = 0.5454545454545454
$ws.Range("N6").Value = 0.7777777777777778
$ws.Range("O6").Value = 0.7654320987654321

# Row 11 (MLP, technique 10%) - full metrics refresh + parameter change
$ws.Range("C11").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 0.5
$ws.Range("J11").Value = 0.5833333333333334
$ws.Range("K11").Value = 0.7777777777777778
$ws.Range("L11").Value = 0.4666666666666667
$ws.Range("M11").Value = 0.2727272727272727
$ws.Range("N11").Value = 0.7777777777777778
$ws.Range("O11").Value = 0.673469387755102

# Row 16 (MLP, technique Free) - full metrics refresh + parameter change
$ws.Range("C16").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 0.3
$ws.Range("J16").Value = 0.4166666666666667
$ws.Range("K16").Value = 0.5555555555555556
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.09090909090909091
$ws.Range("N16").Value = 0.5555555555555556
$ws.Range("O16").Value = 0.5454545454545454
